$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H40").Value = 4769.95
$ws.Range("I40").Value = 3560.1428
$ws.Range("J40").Value = 7592.8335
$ws.Range("K40").Value = 3560.1428
$ws.Range("L40").Value = 7592.8335
$ws.Range("M40").Value = -3385.1428
$ws.Range("N40").Value = -7942.8335

$ws.Range("H94").Value = 3831.9
$ws.Range("I94").Value = 3831.9
$ws.Range("K94").Value = 3831.9
$ws.Range("M94").Value = -3380.9

$ws.Range("H100").Value = 4784.5557
$ws.Range("I100").Value = 3612.2
$ws.Range("K100").Value = 3612.2
$ws.Range("M100").Value = -3071.2

$ws.Range("H113").Value = 4251
$ws.Range("I113").Value = 4001.3333
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4001.3333
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -747.3332999999998
$ws.Range("N113").Value = -11508

$ws.Range("H129").Value = 1538.8667
$ws.Range("I129").Value = 935
$ws.Range("J129").Value = 2444.6667
$ws.Range("K129").Value = 2805
$ws.Range("L129").Value = 7334.000100000001
$ws.Range("M129").Value = 2195
$ws.Range("N129").Value = -17334.0001

$ws.Range("H131").Value = 4144.1665
$ws.Range("I131").Value = 4173
$ws.Range("J131").Value = 4000
$ws.Range("K131").Value = 12519
$ws.Range("L131").Value = 12000
$ws.Range("M131").Value = -7479
$ws.Range("N131").Value = -22080

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 7696965.5
$ws.Range("I32").Value = 2255.4
$ws.Range("J32").Value = 33346000
$ws.Range("K32").Value = 2255.4
$ws.Range("L32").Value = 33346000
$ws.Range("M32").Value = -1968.4
$ws.Range("N32").Value = -33346574

$ws.Range("H44").Value = 13162.8
$ws.Range("J44").Value = 13162.8
$ws.Range("L44").Value = 13162.8
$ws.Range("N44").Value = -14138.8

$ws.Range("H97").Value = 930.8
$ws.Range("I97").Value = 529.6667
$ws.Range("J97").Value = 1532.5
$ws.Range("K97").Value = 529.6667
$ws.Range("L97").Value = 1532.5
$ws.Range("M97").Value = -33.66669999999999
$ws.Range("N97").Value = -2524.5

$ws.Range("H112").Value = 17500.5
$ws.Range("J112").Value = 17500.5
$ws.Range("L112").Value = 17500.5
$ws.Range("N112").Value = -20454.5

$ws.Range("H122").Value = 2502.0908
$ws.Range("I122").Value = 2302.3
$ws.Range("K122").Value = 6906.900000000001
$ws.Range("M122").Value = -4456.900000000001

$ws.Range("H135").Value = 95664
$ws.Range("J135").Value = 95664
$ws.Range("L135").Value = 95664
$ws.Range("N135").Value = -105804

$ws = $wb.Worksheets("BSM")
$ws.Range("H95").Value = 5418.5
$ws.Range("J95").Value = 5418.5
$ws.Range("L95").Value = 5418.5
$ws.Range("N95").Value = -10910.5

$ws.Range("H107").Value = 6083.615
$ws.Range("I107").Value = 1577.6
$ws.Range("J107").Value = 8899.875
$ws.Range("K107").Value = 1577.6
$ws.Range("L107").Value = 8899.875
$ws.Range("M107").Value = 342.4000000000001
$ws.Range("N107").Value = -12739.875

$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180

$ws.Range("H134").Value = 4000
$ws.Range("I134").Value = 4000
$ws.Range("K134").Value = 12000
$ws.Range("M134").Value = -9465

$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 1848.7142
$ws.Range("I22").Value = 592.44446
$ws.Range("J22").Value = 4110
$ws.Range("K22").Value = 592.44446
$ws.Range("L22").Value = 4110
$ws.Range("M22").Value = -242.44446
$ws.Range("N22").Value = -4810

$ws.Range("H62").Value = 3189.889
$ws.Range("J62").Value = 3379.8
$ws.Range("L62").Value = 3379.8
$ws.Range("N62").Value = -4627.8

$ws.Range("H65").Value = 3189.889
$ws.Range("J65").Value = 3379.8
$ws.Range("L65").Value = 16899
$ws.Range("N65").Value = -23139

$ws.Range("H132").Value = 3105.2104
$ws.Range("J132").Value = 4128.5713
$ws.Range("L132").Value = 12385.7139
$ws.Range("N132").Value = -17445.7139

$ws.Range("H141").Value = 109141.164
$ws.Range("J141").Value = 109141.164
$ws.Range("L141").Value = 109141.164
$ws.Range("N141").Value = -119501.164

$ws = $wb.Worksheets("CUL")
$ws.Range("H12").Value = 275.7
$ws.Range("J12").Value = 298
$ws.Range("L12").Value = 894
$ws.Range("N12").Value = -1240

$ws.Range("H107").Value = 527.1875
$ws.Range("I107").Value = 577
$ws.Range("K107").Value = 1731
$ws.Range("M107").Value = 189

$ws.Range("H125").Value = 7989.5
$ws.Range("I125").Value = 7989.5
$ws.Range("K125").Value = 23968.5
$ws.Range("M125").Value = -19048.5

$ws.Range("H136").Value = 11000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws.Range("H139").Value = 3674.2727
$ws.Range("I139").Value = 3271.3333
$ws.Range("K139").Value = 9813.999899999999
$ws.Range("M139").Value = -4673.999899999999

$ws.Range("H141").Value = 750
$ws.Range("I141").Value = 750
$ws.Range("K141").Value = 2250
$ws.Range("M141").Value = 2930

$ws = $wb.Worksheets("GSM")
$ws.Range("H80").Value = 1909.3
$ws.Range("J80").Value = 1598.3334
$ws.Range("L80").Value = 1598.3334
$ws.Range("N80").Value = -3594.3334

$ws.Range("H83").Value = 1909.3
$ws.Range("J83").Value = 1598.3334
$ws.Range("L83").Value = 7991.666999999999
$ws.Range("N83").Value = -17975.667

$ws.Range("H102").Value = 1739.5454
$ws.Range("I102").Value = 1739.5454
$ws.Range("K102").Value = 1739.5454
$ws.Range("M102").Value = -117.5454

$ws.Range("H122").Value = 2538.2
$ws.Range("I122").Value = 2672.75
$ws.Range("K122").Value = 8018.25
$ws.Range("M122").Value = -5568.25

$ws.Range("H126").Value = 3747.25
$ws.Range("I126").Value = 1744.5
$ws.Range("K126").Value = 5233.5
$ws.Range("M126").Value = -2763.5

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 3341.5881
$ws.Range("I7").Value = 2528.1333
$ws.Range("K7").Value = 2528.1333
$ws.Range("M7").Value = -2416.1333

$ws.Range("H22").Value = 910.4545000000001
$ws.Range("I22").Value = 863.3333
$ws.Range("J22").Value = 1122.5
$ws.Range("K22").Value = 863.3333
$ws.Range("L22").Value = 1122.5
$ws.Range("M22").Value = -568.3333
$ws.Range("N22").Value = -1712.5

$ws.Range("H27").Value = 910.4545000000001
$ws.Range("I27").Value = 863.3333
$ws.Range("J27").Value = 1122.5
$ws.Range("K27").Value = 863.3333
$ws.Range("L27").Value = 1122.5
$ws.Range("M27").Value = -756.3333
$ws.Range("N27").Value = -1336.5

$ws.Range("H55").Value = 505
$ws.Range("I55").Value = 1001
$ws.Range("J55").Value = 381
$ws.Range("K55").Value = 1001
$ws.Range("L55").Value = 381
$ws.Range("M55").Value = -828
$ws.Range("N55").Value = -727

$ws.Range("H97").Value = 27333.334
$ws.Range("J97").Value = 27333.334
$ws.Range("L97").Value = 27333.334
$ws.Range("N97").Value = -29315.334

$ws.Range("H100").Value = 5833.1665
$ws.Range("I100").Value = 4437.25
$ws.Range("J100").Value = 8625
$ws.Range("K100").Value = 4437.25
$ws.Range("L100").Value = 8625
$ws.Range("M100").Value = -3896.25
$ws.Range("N100").Value = -9707

$ws.Range("H122").Value = 4166.6665
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -8800
$ws.Range("N122").Value = -19900

$ws.Range("H126").Value = 3341.5881
$ws.Range("I126").Value = 2528.1333
$ws.Range("K126").Value = 7584.3999
$ws.Range("M126").Value = -5114.3999

$ws.Range("H132").Value = 4007.5
$ws.Range("I132").Value = 1828.7273
$ws.Range("J132").Value = 8800.799999999999
$ws.Range("K132").Value = 5486.1819
$ws.Range("L132").Value = 26402.4
$ws.Range("M132").Value = -2956.1819
$ws.Range("N132").Value = -31462.4

$ws = $wb.Worksheets("WVR")
$ws.Range("H95").Value = 18814.666
$ws.Range("J95").Value = 18814.666
$ws.Range("L95").Value = 18814.666
$ws.Range("N95").Value = -24306.666

$ws.Range("H107").Value = 767.3077
$ws.Range("I107").Value = 484.5
$ws.Range("K107").Value = 1453.5
$ws.Range("M107").Value = 466.5

$ws.Range("H126").Value = 5690.5713
$ws.Range("I126").Value = 3958.5
$ws.Range("K126").Value = 11875.5
$ws.Range("M126").Value = -9405.5

$ws.Range("H132").Value = 1178.2858
$ws.Range("I132").Value = 1178.2858
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3534.8574
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1004.8574
$ws.Range("N132").ClearContents()
